$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) is stored as text so values like "1.00" or "63.365.26"
# are not coerced into numbers and lose formatting / precision.
$ws.Range("D2:D51").NumberFormat = "@"

$updates = @{
    'D2' = '63.365.26'
    'E2' = '  -5.28%  '
    'D3' = '2.961.95'
    'E3' = '  -8.00%  '
    'D4' = '1.00'
    'E4' = '  +0.03%  '
    'D5' = '547.77'
    'E5' = '  -5.55%  '
    'D6' = '134.93'
    'E6' = '  -10.86%  '
    'D7' = '0.998'
    'E7' = '  -0.03%  '
    'D8' = '2.955.05'
    'E8' = '  -7.98%  '
    'D9' = '0.472'
    'E9' = '  -13.02%  '
    'E10' = '  -13.93%  '
    'D11' = '5.93'
    'E11' = '  -13.23%  '
    'D12' = '0.448'
    'E12' = '  -10.51%  '
    'D13' = '33.73'
    'E13' = '  -11.99%  '
    'E14' = '  -12.80%  '
    'D15' = '3.422.38'
    'E15' = '  -8.50%  '
    'D16' = '63.240.93'
    'E16' = '  -5.35%  '
    'D17' = '0.109'
    'E17' = '  -4.18%  '
    'D18' = '2.950.13'
    'E18' = '  -8.35%  '
    'D19' = '6.39'
    'E19' = '  -10.39%  '
    'D20' = '462.23'
    'E20' = '  -14.82%  '
    'E21' = '  -13.27%  '
    'D22' = '0.647'
    'E22' = '  -14.49%  '
    'D23' = '6.81'
    'E23' = '  -11.87%  '
    'B24' = 'Litecoin'
    'C24' = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    'D24' = '75.39'
    'E24' = '  -11.91%  '
    'B25' = 'InternetComputer(DFINITY)'
    'C25' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D25' = '11.97'
    'E25' = '  -11.02%  '
    'E26' = '  +0.06%  '
    'D27' = '2.63'
    'E27' = '  -16.45%  '
    'D28' = '2.00'
    'E28' = '  -5.39%  '
    'E29' = '  -0.37%  '
    'E30' = '  -8.66%  '
    'D31' = '25.10'
    'E31' = '  -14.38%  '
    'D32' = '2.53'
    'E32' = '  -2.03%  '
    'D33' = '1.06'
    'E33' = '  -6.58%  '
    'D34' = '486.20'
    'E34' = '  -10.57%  '
    'B35' = 'NEARProtocol'
    'C35' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D35' = '5.10'
    'E35' = '  -10.62%  '
    'B36' = 'OKB'
    'C36' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'D36' = '51.76'
    'E36' = '  -2.52%  '
    'D37' = '5.62'
    'E37' = '  -14.29%  '
    'D38' = '0.0386'
    'E38' = '  -10.96%  '
    'D39' = '0.0754'
    'E39' = '  -10.91%  '
    'E40' = '  -7.39%  '
    'D41' = '8.02'
    'E41' = '  -12.46%  '
    'D42' = '2.804.83'
    'E42' = '  -4.42%  '
    'B43' = 'USDe'
    'C43' = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
    'D43' = '0.998'
    'E43' = '  -0.23%  '
    'B44' = 'dogwifhat'
    'C44' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D44' = '2.35'
    'E44' = '  -8.64%  '
    'E45' = '  -12.50%  '
    'B46' = 'Fetch.AI'
    'C46' = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    'D46' = '1.94'
    'E46' = '  -8.13%  '
    'B47' = 'Monero'
    'C47' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D47' = '115.92'
    'E47' = '  -6.15%  '
    'B48' = 'InjectiveProtocol'
    'C48' = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    'D48' = '23.68'
    'E48' = '  -8.74%  '
    'B49' = 'PEPE'
    'C49' = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
    'D49' = '0.0₃0504'
    'E49' = '  -13.84%  '
    'B50' = 'Stellar'
    'C50' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D50' = '0.104'
    'E50' = '  -8.22%  '
    'D51' = '1.93'
    'E51' = '  -18.75%  '
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$wb.Save()
